$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits (wording updates to the plan/timeline activities) ---
# Week 1 activity: softened wording ("Start with" instead of "I'd recommend starting with")
$ws.Range("B2").Value = "Undertake wider data exploration and reading to find a topic you might be interested in progressing with. Start with a data search first to avoid proposing a research area then forcing potentially unsuitable data to fit it."

# Week 5 activity: replaced "Complete a draft of part 1..." with new methodology-draft wording
$ws.Range("B6").Value = "Create a draft of the methodology / initial ideas that can be expanded on as we move to more advanced analysis in the second part of the course."

# "Before submission" activity: reworded to reference introduction/lit review instead of "part 1"
$ws.Range("B12").Value = "Check that your assignment follows the standard model of scientific investigation. You may need to update your introduction and literature review based on the analysis you undertook to ensure a seamless narrative throughout.  "

# --- Style changes: drop the cell borders for B2, B11, B12 (was bordered style, now borderless) ---
# NOTE: Borders must be cleared per-cell; clearing across a multi-cell range only affects the first cell.
foreach ($addr in @("B2", "B11", "B12")) {
    $rng = $ws.Range($addr)
    $rng.Borders.Item(7).LineStyle = 0
    $rng.Borders.Item(8).LineStyle = 0
    $rng.Borders.Item(9).LineStyle = 0
    $rng.Borders.Item(10).LineStyle = 0
}

# --- Style changes: drop borders for B8:B10 and justify the text ---
foreach ($addr in @("B8", "B9", "B10")) {
    $rng = $ws.Range($addr)
    $rng.Borders.Item(7).LineStyle = 0
    $rng.Borders.Item(8).LineStyle = 0
    $rng.Borders.Item(9).LineStyle = 0
    $rng.Borders.Item(10).LineStyle = 0
}
$ws.Range("B8:B10").HorizontalAlignment = -4130

# --- Row height: week-5 row shrinks to match its shorter new text ---
$ws.Rows.Item(6).RowHeight = 30.25

# --- View state: the author's selection ended up on B12 ---
[void]$ws.Range("B12").Select()
